$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MapDesign_lvl1")

# Paint the floor/wall cells: yellow fill (new style) + value 5
$floorRanges = @(
    "F2:AD2",
    "F3",
    "AD3",
    "F4",
    "AD4",
    "F5",
    "AD5",
    "F6",
    "AD6",
    "F7",
    "AD7",
    "F8",
    "AD8",
    "F9",
    "M9",
    "T9",
    "W9",
    "AD9",
    "F10",
    "M10",
    "T10",
    "W10",
    "AD10",
    "F11",
    "M11",
    "T11",
    "W11",
    "AA11",
    "AD11",
    "F12",
    "M12",
    "T12",
    "W12",
    "AA12",
    "AD12",
    "F13",
    "M13",
    "T13",
    "W13",
    "AA13",
    "AD13",
    "F14",
    "M14",
    "T14",
    "W14",
    "AA14",
    "AD14",
    "F15",
    "M15",
    "T15",
    "W15",
    "AA15",
    "AD15",
    "F16",
    "M16",
    "T16",
    "W16",
    "AA16",
    "AD16",
    "F17",
    "M17",
    "T17",
    "AA17",
    "AD17",
    "F18",
    "N18:S18",
    "AA18",
    "AD18",
    "F19",
    "AA19",
    "AD19",
    "F20",
    "AD20",
    "F21",
    "AD21",
    "F22",
    "AD22",
    "F23",
    "AD23"
)

foreach ($addr in $floorRanges) {
    $rng = $ws.Range($addr)
    $rng.Interior.Color = 65535
    $rng.Value = 5
}

# Update the active selection to match the new edit focus
[void]$ws.Range("AA11:AA19").Select()

